$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $donorRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

$ws.Range("E2").Value = "2026-02-23 18:49:09"
Set-TextValue "H2" "C2" "35%"
$ws.Range("E3").Value = "2026-02-23 18:49:12"
$ws.Range("E4").Value = "2026-02-23 18:49:15"
$ws.Range("J4").Value = "1024.7 hPa"
$ws.Range("E5").Value = "2026-02-23 18:49:17"
$ws.Range("E6").Value = "2026-02-23 18:49:20"
$ws.Range("E7").Value = "2026-02-23 18:49:23"
$ws.Range("E8").Value = "2026-02-23 18:49:25"
$ws.Range("K8").Value = "15.8 MJ/m2"
$ws.Range("E9").Value = "2026-02-23 18:49:28"
Set-TextValue "H9" "C9" "71%"
$ws.Range("E10").Value = "2026-02-23 18:49:31"
$ws.Range("E11").Value = "2026-02-23 18:49:33"
$ws.Range("E12").Value = "2026-02-23 18:49:36"
Set-TextValue "H12" "C12" "83%"
$ws.Range("E13").Value = "2026-02-23 18:49:38"
Set-TextValue "H13" "C13" "58%"
$ws.Range("J13").Value = "1026.8 hPa"
$ws.Range("O13").Value = "7.2 °C"
$ws.Range("E14").Value = "2026-02-23 18:49:41"
Set-TextValue "H14" "C14" "72%"
$ws.Range("E15").Value = "2026-02-23 18:49:43"
$ws.Range("O15").Value = "12.9 °C"
$ws.Range("E16").Value = "2026-02-23 18:49:46"
$ws.Range("E17").Value = "2026-02-23 18:49:48"
$ws.Range("E18").Value = "2026-02-23 18:49:51"
$ws.Range("O18").Value = "11.2 °C"
$ws.Range("E19").Value = "2026-02-23 18:49:54"
Set-TextValue "H19" "C19" "44%"
$ws.Range("E20").Value = "2026-02-23 18:49:56"
$ws.Range("K20").Value = "16.6 MJ/m2"
$ws.Range("E21").Value = "2026-02-23 18:49:59"
$ws.Range("J21").Value = "1025.8 hPa"
$ws.Range("O21").Value = "9.7 °C"
$ws.Range("E22").Value = "2026-02-23 18:50:02"
Set-TextValue "H22" "C22" "22%"
$ws.Range("O22").Value = "3.4 °C"
$ws.Range("E23").Value = "2026-02-23 18:50:05"
$ws.Range("O23").Value = "3.8 °C"
$ws.Range("E24").Value = "2026-02-23 18:50:07"
Set-TextValue "H24" "C24" "81%"
$ws.Range("J24").Value = "1026.2 hPa"
$ws.Range("E25").Value = "2026-02-23 18:50:10"
Set-TextValue "H25" "C25" "26%"
$ws.Range("E26").Value = "2026-02-23 18:50:12"
Set-TextValue "H26" "C26" "49%"
$ws.Range("E27").Value = "2026-02-23 18:50:15"
$ws.Range("E28").Value = "2026-02-23 18:50:18"
$ws.Range("O28").Value = "11.3 °C"
$ws.Range("E29").Value = "2026-02-23 18:50:20"
Set-TextValue "H29" "C29" "81%"
$ws.Range("E30").Value = "2026-02-23 18:50:23"
Set-TextValue "H30" "C30" "68%"
$ws.Range("E31").Value = "2026-02-23 18:50:26"
$ws.Range("J31").Value = "1023.9 hPa"
$ws.Range("E32").Value = "2026-02-23 18:50:28"
$ws.Range("O32").Value = "8.4 °C"
$ws.Range("E33").Value = "2026-02-23 18:50:31"
$ws.Range("J33").Value = "1025.3 hPa"
$ws.Range("O33").Value = "8.8 °C"
$ws.Range("E34").Value = "2026-02-23 18:50:34"
Set-TextValue "H34" "C34" "40%"
$ws.Range("E35").Value = "2026-02-23 18:50:36"
Set-TextValue "H35" "C35" "37%"
$ws.Range("J35").Value = "1025.0 hPa"
$ws.Range("O35").Value = "12.8 °C"
$ws.Range("E36").Value = "2026-02-23 18:50:39"
Set-TextValue "H36" "C36" "72%"
$ws.Range("J36").Value = "1024.7 hPa"
$ws.Range("E37").Value = "2026-02-23 18:50:42"
Set-TextValue "H37" "C37" "64%"
$ws.Range("O37").Value = "9.5 °C"
$ws.Range("E38").Value = "2026-02-23 18:50:44"
$ws.Range("E39").Value = "2026-02-23 18:50:47"
$ws.Range("E40").Value = "2026-02-23 18:50:50"
Set-TextValue "H40" "C40" "60%"
$ws.Range("O40").Value = "9.1 °C"
$ws.Range("E41").Value = "2026-02-23 18:50:52"
Set-TextValue "H41" "C41" "70%"
$ws.Range("K41").Value = "15.8 MJ/m2"
$ws.Range("E42").Value = "2026-02-23 18:50:55"
Set-TextValue "H42" "C42" "77%"
$ws.Range("O42").Value = "12.0 °C"
$ws.Range("E43").Value = "2026-02-23 18:50:57"
$ws.Range("O43").Value = "10.4 °C"
$ws.Range("E44").Value = "2026-02-23 18:51:00"
Set-TextValue "H44" "C44" "32%"
$ws.Range("O44").Value = "3.4 °C"
$ws.Range("E45").Value = "2026-02-23 18:51:03"
Set-TextValue "H45" "C45" "50%"
$ws.Range("O45").Value = "8.7 °C"
$ws.Range("E46").Value = "2026-02-23 18:51:05"
$ws.Range("J46").Value = "1026.1 hPa"
$ws.Range("O46").Value = "10.4 °C"
